# Update "horarios" (schedules) workbook with a fresh scrape snapshot.
# New scrape time: 03:19:42 (previously 02:54:27)

$wb = $excel.ActiveWorkbook

$oldTime = "02:54:27"
$newTime = "03:19:42"

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"

$ws1.Range("A6").Value  = $newTime
$ws1.Range("B6").Value  = "03:48"
$ws1.Range("C6").Value  = "14_ABASTO"
$ws1.Range("D6").Value  = 29

$ws1.Range("A7").Value  = $newTime
$ws1.Range("B7").Value  = "04:01"
$ws1.Range("C7").Value  = "81_EL PELIGRO"
$ws1.Range("D7").Value  = 42

$ws1.Range("A8").Value  = $newTime
$ws1.Range("B8").Value  = "04:46"
$ws1.Range("C8").Value  = "215_EL PELIGRO"
$ws1.Range("D8").Value  = 87

$ws1.Range("A9").Value  = $newTime
$ws1.Range("B9").Value  = "04:53"
$ws1.Range("C9").Value  = "11_ETCHEVERRY"
$ws1.Range("D9").Value  = 94

$ws1.Range("A10").Value = $newTime
$ws1.Range("B10").Value = "05:11"
$ws1.Range("C10").Value = "17_ROMERO"
$ws1.Range("D10").Value = 112

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 87

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
